$wb = $excel.ActiveWorkbook

# --- Rename the worksheets ---
$wsCompare = $wb.Worksheets.Item(1)
$wsCompare.Name = "Feature_Method_Compare"

$wsAllTeams = $wb.Worksheets.Item(2)
$wsAllTeams.Name = "AllTeams_PrimaryFeatures"

# --- Feature_Method_Compare sheet data (rows 2-15) ---
$compareData = @(
    @("Random Forest Importance", "Random Forest", 0.034209786, 0.759611941),
    @("Random Forest Importance", "Decision Tree", $null, $null),
    @("Pearson Correlation", "Random Forest", 0.266958118, -0.87588263),
    @("Pearson Correlation", "Decision Tree", $null, $null),
    @("Spearman Correlation", "Random Forest", 0.266958118, -0.87588263),
    @("Spearman Correlation", "Decision Tree", $null, $null),
    @("Kendall Correlation", "Random Forest", 0.266958118, -0.87588263),
    @("Kendall Correlation", "Decision Tree", $null, $null),
    @("RFECV", "Random Forest", 0.028030902, 0.803030219),
    @("RFECV", "Decision Tree", $null, $null),
    @("Engineered Features", "Random Forest", 0.158268865, -0.112136306),
    @("Engineered Features", "Decision Tree", $null, $null),
    @("OPS+ERA (Control)", "Random Forest", 0.266958118, -0.87588263),
    @("OPS+ERA (Control)", "Decision Tree", $null, $null),
)

$r = 2
foreach ($row in $compareData) {
    $wsCompare.Cells.Item($r, 1).Value = $row[0]
    $wsCompare.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -eq $null) {
        $wsCompare.Cells.Item($r, 3).Value = ""
    } else {
        $wsCompare.Cells.Item($r, 3).Value = $row[2]
    }
    if ($row[3] -eq $null) {
        $wsCompare.Cells.Item($r, 4).Value = ""
    } else {
        $wsCompare.Cells.Item($r, 4).Value = $row[3]
    }
    $r++
}

# --- AllTeams_PrimaryFeatures sheet data (rows 2-65) ---
$teamsData = @(
    @("ARI", "Random Forest", 0.311020915, 0.021888774),
    @("ARI", "Decision Tree", $null, $null),
    @("ATL", "Random Forest", 0.465789795, 0.639750301),
    @("ATL", "Decision Tree", $null, $null),
    @("BAL", "Random Forest", 1.349911116, -0.243245939),
    @("BAL", "Decision Tree", $null, $null),
    @("BOS", "Random Forest", 1.040190748, -0.376067247),
    @("BOS", "Decision Tree", $null, $null),
    @("CHC", "Random Forest", 0.545757453, -0.32751522),
    @("CHC", "Decision Tree", $null, $null),
    @("CIN", "Random Forest", 0.748108521, -1.65958995),
    @("CIN", "Decision Tree", $null, $null),
    @("CLE", "Random Forest", 0.408579042, 0.062406277),
    @("CLE", "Decision Tree", $null, $null),
    @("COL", "Random Forest", 0.187925682, 0.26865064),
    @("COL", "Decision Tree", $null, $null),
    @("CWS", "Random Forest", 1.677013769, -0.585264382),
    @("CWS", "Decision Tree", $null, $null),
    @("DET", "Random Forest", 0.813812561, -0.710199546),
    @("DET", "Decision Tree", $null, $null),
    @("HOU", "Random Forest", 1.461536637, -10.001451313),
    @("HOU", "Decision Tree", $null, $null),
    @("KC", "Random Forest", 0.5904941379999999, 0.006763876),
    @("KC", "Decision Tree", $null, $null),
    @("LAA", "Random Forest", 0.123569487, -6.944010577),
    @("LAA", "Decision Tree", $null, $null),
    @("LAD", "Random Forest", 0.540866382, 0.121571951),
    @("LAD", "Decision Tree", $null, $null),
    @("MIA", "Random Forest", 1.037784616, -41.65506127),
    @("MIA", "Decision Tree", $null, $null),
    @("MIL", "Random Forest", 0.51705843, -0.418231115),
    @("MIL", "Decision Tree", $null, $null),
    @("MIN", "Random Forest", 1.650620164, -0.108594292),
    @("MIN", "Decision Tree", $null, $null),
    @("NYM", "Random Forest", 0.694597192, -1.39191152),
    @("NYM", "Decision Tree", $null, $null),
    @("NYY", "Random Forest", 0.674551363, -12.238759491),
    @("NYY", "Decision Tree", $null, $null),
    @("OAK", "Random Forest", 1.9373596, 0.09381186699999999),
    @("OAK", "Decision Tree", $null, $null),
    @("PHI", "Random Forest", 0.231009424, 0.387260063),
    @("PHI", "Decision Tree", $null, $null),
    @("PIT", "Random Forest", 0.658129657, 0.226797474),
    @("PIT", "Decision Tree", $null, $null),
    @("SD", "Random Forest", 1.103639356, -0.134704441),
    @("SD", "Decision Tree", $null, $null),
    @("SEA", "Random Forest", 0.259394609, -0.016585307),
    @("SEA", "Decision Tree", $null, $null),
    @("SF", "Random Forest", 0.395505102, -3.455050431),
    @("SF", "Decision Tree", $null, $null),
    @("STL", "Random Forest", 0.784943619, -2.130896984),
    @("STL", "Decision Tree", $null, $null),
    @("TB", "Random Forest", 1.517802631, 0.09236860400000001),
    @("TB", "Decision Tree", $null, $null),
    @("TEX", "Random Forest", 1.91439705, -0.361941598),
    @("TEX", "Decision Tree", $null, $null),
    @("TOR", "Random Forest", 0.264794774, -30.120306641),
    @("TOR", "Decision Tree", $null, $null),
    @("WAS", "Random Forest", 0.237297393, 0.685054775),
    @("WAS", "Decision Tree", $null, $null),
    @("All", "Random Forest", 0.034209786, 0.759611941),
    @("All", "Decision Tree", $null, $null),
    @("MLB_years", "Random Forest", 0.000107512, -180.950046875),
    @("MLB_years", "Decision Tree", $null, $null),
)

$r = 2
foreach ($row in $teamsData) {
    $wsAllTeams.Cells.Item($r, 1).Value = $row[0]
    $wsAllTeams.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -eq $null) {
        $wsAllTeams.Cells.Item($r, 3).Value = ""
    } else {
        $wsAllTeams.Cells.Item($r, 3).Value = $row[2]
    }
    if ($row[3] -eq $null) {
        $wsAllTeams.Cells.Item($r, 4).Value = ""
    } else {
        $wsAllTeams.Cells.Item($r, 4).Value = $row[3]
    }
    $r++
}
